$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.654
$ws.Range("C5").Value = -12.156
$ws.Range("D7").Value = -7.325999999999999
$ws.Range("C9").Value = -11.671
$ws.Range("C11").Value = -12.994
$ws.Range("D11").Value = -8.113
$ws.Range("A21").Value = -21.547
$ws.Range("C21").Value = -12.679
$ws.Range("D21").Value = -7.965999999999999
$ws.Range("A23").Value = -21.587
$ws.Range("A25").Value = -21.937
